$d = $word.ActiveDocument

$replacements = @(
    @("42×85=3570", "17×42=714"),
    @("47×84=3948", "57×93=5301"),
    @("40×51=2040", "49×61=2989"),
    @("50×14=700", "47×28=1316"),
    @("48×65=3120", "63×91=5733"),
    @("33×12=396", "66×93=6138"),
    @("37×62=2294", "12×14=168"),
    @("80×73=5840", "34×52=1768"),
    @("75×22=1650", "50×79=3950"),
    @("72×40=2880", "21×32=672"),
    @("25×94=2350", "92×55=5060"),
    @("46×30=1380", "40×43=1720"),
    @("63×36=2268", "30×17=510"),
    @("61×72=4392", "66×79=5214"),
    @("88×84=7392", "95×52=4940"),
    @("96×73=7008", "94×85=7990"),
    @("63×68=4284", "69×11=759"),
    @("47×93=4371", "62×76=4712"),
    @("87×32=2784", "67×95=6365"),
    @("53×86=4558", "83×50=4150"),
    @("19×62=1178", "48×96=4608"),
    @("66×19=1254", "15×16=240"),
    @("30×92=2760", "26×78=2028"),
    @("15×64=960", "43×78=3354"),
    @("93×84=7812", "84×46=3864")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
